$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Add the new row of data (row 9) below the existing table.
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Rettifica indirizzo post accertamento"

# Match formatting used by the rest of the table (font + thin border, center/wrap).
$ws.Range("A9").Style = $ws.Range("A2").Style
$ws.Range("B9").Style = $ws.Range("B2").Style

$ws.Range("A9:B9").Borders.LineStyle = 1
$ws.Range("A9:B9").Borders.Weight = 2
$ws.Range("A9:B9").Borders.ColorIndex = 64

# Restore the active cell/selection used when the file was last saved.
$ws.Range("B13").Select()
